$wb = $excel.ActiveWorkbook

# --- Folder Inventory sheet: swap rows 3 and 4 (columns A and B) ---
$ws1 = $wb.Worksheets.Item("Folder Inventory")

$ws1.Range("A3").Value = "Azure Virtual Machine And Compute"
$ws1.Range("B3").Value = "Azure Virtual Machine And Compute"

$ws1.Range("A4").Value = "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals"
$ws1.Range("B4").Value = "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals"

# --- Metadata sheet updates ---
$ws2 = $wb.Worksheets.Item("Metadata")

$ws2.Range("B2").Value = "CloudLabsAI-Azure/MS-Innovation-Release-Notes"
$ws2.Range("B3").Value = "2025-06-12 10:36:39 UTC"

# B5 ("Workflow Run") is stored as TEXT ("2" -> "10"), not a number.
# Force text type via NumberFormat while writing, then restore the
# format so no stray style is left behind on the cell.
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "10"
$ws2.Range("B5").Style = "Normal"
